$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 7 fresh rows right after the current last data row (22) and
#    before the signature rows (27:28). This naturally pushes the signature
#    rows down to 34:35 while preserving the blank gap, and fixes up the
#    merged-cell references for those rows automatically.
# ---------------------------------------------------------------------------
$ws.Rows("23:29").Insert()

# ---------------------------------------------------------------------------
# 2) Copy the formatting of the existing data rows onto the new rows so the
#    borders/styles match: rows 16-21 use the "inner" row style, the old
#    row 22 used the heavier "closing" row style (now needed at row 29).
#    Row 22 itself is no longer the closing row of the table, so restyle it
#    to the "inner" row look first (before its old formatting is copied to
#    row 29 below - order matters).
# ---------------------------------------------------------------------------
$ws.Range("B16:J21").Copy()
$ws.Range("B23:J28").PasteSpecial(-4122)

$ws.Range("B22:J22").Copy()
$ws.Range("B29:J29").PasteSpecial(-4122)

$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Header / label text (content unchanged, only shared-string slot shifts
#    under the hood - writing the literal text keeps it correct either way).
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "ESTADO DE CUENTA"
$ws.Range("B7").Value = "RAZON SOCIAL:"
$ws.Range("B11").Value = "VALOR MORA"
$ws.Range("B13").Value = "Cant. Trabajadores"
$ws.Range("E13").Value = "Cant. Periodos"
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"
$ws.Range("J15").Value = "Observaciones"

# ---------------------------------------------------------------------------
# 4) Updated summary figures.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 2505432
$ws.Range("C13").Value = 2

# ---------------------------------------------------------------------------
# 5) Existing worker (DAYANA) rows 16-22: periods now run newest -> oldest
#    (2305 .. 2211) and the "155369" balance moved off row 16 onto row 22.
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "45530730"
$ws.Range("D16").Value = "DAYANA ROSSY GONZALEZ STAND"
$ws.Range("E16").Value = "2305"
$ws.Range("F16").Value = 124295
$ws.Range("G16").Value = 3884234

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45530730"
$ws.Range("D17").Value = "DAYANA ROSSY GONZALEZ STAND"
$ws.Range("E17").Value = "2304"
$ws.Range("F17").Value = 155369
$ws.Range("G17").Value = 3884234

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "45530730"
$ws.Range("D18").Value = "DAYANA ROSSY GONZALEZ STAND"
$ws.Range("E18").Value = "2303"
$ws.Range("F18").Value = 155369
$ws.Range("G18").Value = 3884234

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45530730"
$ws.Range("D19").Value = "DAYANA ROSSY GONZALEZ STAND"
$ws.Range("E19").Value = "2302"
$ws.Range("F19").Value = 155369
$ws.Range("G19").Value = 3884234

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "45530730"
$ws.Range("D20").Value = "DAYANA ROSSY GONZALEZ STAND"
$ws.Range("E20").Value = "2301"
$ws.Range("F20").Value = 155369
$ws.Range("G20").Value = 3884234

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "45530730"
$ws.Range("D21").Value = "DAYANA ROSSY GONZALEZ STAND"
$ws.Range("E21").Value = "2212"
$ws.Range("F21").Value = 155369
$ws.Range("G21").Value = 3884234

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "45530730"
$ws.Range("D22").Value = "DAYANA ROSSY GONZALEZ STAND"
$ws.Range("E22").Value = "2211"
$ws.Range("F22").Value = 155369
$ws.Range("G22").Value = 3884234

# ---------------------------------------------------------------------------
# 6) New worker (YESENIA) rows 23-29, periods 2305 .. 2211.
# ---------------------------------------------------------------------------
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1047413615"
$ws.Range("D23").Value = "YESENIA PAOLA ORTEGA MARTINEZ"
$ws.Range("E23").Value = "2305"
$ws.Range("F23").Value = 170461
$ws.Range("G23").Value = 5326934

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1047413615"
$ws.Range("D24").Value = "YESENIA PAOLA ORTEGA MARTINEZ"
$ws.Range("E24").Value = "2304"
$ws.Range("F24").Value = 213077
$ws.Range("G24").Value = 5326934

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1047413615"
$ws.Range("D25").Value = "YESENIA PAOLA ORTEGA MARTINEZ"
$ws.Range("E25").Value = "2303"
$ws.Range("F25").Value = 213077
$ws.Range("G25").Value = 5326934

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1047413615"
$ws.Range("D26").Value = "YESENIA PAOLA ORTEGA MARTINEZ"
$ws.Range("E26").Value = "2302"
$ws.Range("F26").Value = 213077
$ws.Range("G26").Value = 5326934

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1047413615"
$ws.Range("D27").Value = "YESENIA PAOLA ORTEGA MARTINEZ"
$ws.Range("E27").Value = "2301"
$ws.Range("F27").Value = 213077
$ws.Range("G27").Value = 5326934

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1047413615"
$ws.Range("D28").Value = "YESENIA PAOLA ORTEGA MARTINEZ"
$ws.Range("E28").Value = "2212"
$ws.Range("F28").Value = 213077
$ws.Range("G28").Value = 5326934

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1047413615"
$ws.Range("D29").Value = "YESENIA PAOLA ORTEGA MARTINEZ"
$ws.Range("E29").Value = "2211"
$ws.Range("F29").Value = 213077
$ws.Range("G29").Value = 5326934

# ---------------------------------------------------------------------------
# 7) Signature rows, now at 34:35 (text content unchanged).
# ---------------------------------------------------------------------------
$ws.Range("B34").Value = "___________________________________"
$ws.Range("H34").Value = "___________________________________"
$ws.Range("B35").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H35").Value = "FIRMA DEL REPRESENTANTE LEGAL"

# ---------------------------------------------------------------------------
# 8) Column D widened to fit the new, longer worker name (target ~33.27).
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 32.45
